# Automatic hashcode update: replace specific hashcode values in column B
# of the active worksheet (data/metadata/hashcode.xlsx), matching the
# author commit "Actualizacion automatica hashcode".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B11").Value = "1352d9b99bf06626ff80952eda02d7d2"  # was 1f682c4baf00039722b9d3b2a8f6431f
$ws.Range("B34").Value = "c61e0c5fa0c3d3aeb7f195c62229f494"  # was 9b5fa738b68a8c46f512c3e8ae609d3b
$ws.Range("B44").Value = "a2cfcbfef9b7b4aed5ed06cdf76e820f"  # was 775da89266fde57dfe7ca7c89abf5d91
$ws.Range("B74").Value = "9555bf74da8a390313ded720eb47dce7"  # was 8a74666dc4ebb183229cedc771aa374f
$ws.Range("B89").Value = "160ee88f449d69ffbf488ebe9d2dcc44"  # was e5a9c26e094a5557ae9c4aa83e416d55
$ws.Range("B99").Value = "ec5bd2a050b8a245967e920be6cdaaa2"  # was 0c473cacc596f7b80f753639d0d0ca9c
$ws.Range("B110").Value = "4050bd447a74401c61ea746f9711d4fc"  # was 8c9098805d070995ea6995c660cc73a1
$ws.Range("B121").Value = "27c1bb70cb640d5ca20a759347c927c8"  # was 81667d4f5140992663fc6287a415e11f
$ws.Range("B154").Value = "e9828e955ed4896624069e2230da5da2"  # was 0164192226833e8b2508d9634b0ba903
$ws.Range("B160").Value = "f3de5288eeaf606f566c40f38f1f948a"  # was adf3c1215f1ec05392a34e4fcab6d818
$ws.Range("B161").Value = "9bb4c7968671c6ffbee5b3db18131f17"  # was 1e5c3f3bf56fea72588394470e1cc359
$ws.Range("B162").Value = "28b7081ddd8b2bf574091a34d8703cef"  # was 537a5222143850acb0b8e7c2a56d1a6f
$ws.Range("B168").Value = "36c8cd53ba8a46717318adc0a51706b1"  # was bc95cae257a5ff8399d8aa38ac0096e0
$ws.Range("B180").Value = "4452182d4a3e39871668d09fdb6c1e5b"  # was 8e3e66726412138b9c21d57bc4009d98
$ws.Range("B191").Value = "c73e5ad0a567948972aa3db3a087d497"  # was aec159b771e496e8cb54e48f8a239e8e
$ws.Range("B213").Value = "e11742ebab986b101aaf472dd8371e81"  # was f1a3da6a4991d211f4d0e18b9486ed7a
$ws.Range("B278").Value = "4f4e6e1d7f91885a3a4f184b8ac396e3"  # was 9283cf6e227051ed64790cd8214746ac
$ws.Range("B293").Value = "21201fdc44ce87e98d9209da669acf6b"  # was a7d0b31354aa502f18e0103883abbc31
$ws.Range("B335").Value = "ecbe729ac86df7acbe5e7934836f2f14"  # was fa67257d9e82773e7b9d6f5b58515c14
$ws.Range("B345").Value = "183913fecc02620ae6913e0667b17656"  # was 3d3502f758d76be92c0f4e2ea3201dd1
$ws.Range("B461").Value = "b11b80ec3b93464d6b97a5f9c1948435"  # was 060072cb4a449d58d07838c00b609f70
$ws.Range("B480").Value = "f23b3dca7b162c63f81a3379142179f4"  # was 1fd9ef0f8869fc52d6c81138b24ec41c
$ws.Range("B506").Value = "51d94fbb108c060af0774f3dfc25fd2e"  # was aa1791820592e49d2dde3aff5748084a
$ws.Range("B514").Value = "1ff4dd27e25e4cecffa8c888a063c5c2"  # was 0163ad4ebad868ebcb1fb1d515410e6b
$ws.Range("B524").Value = "586802b4d9ba45de50d961c63708f3c0"  # was b8463e643f40c14c051b7aa3e19cc647
$ws.Range("B534").Value = "76da3783aa2a61aa6867b6ba825b3179"  # was b4d216af1c0225064ccc574065e16246
$ws.Range("B547").Value = "12134a6651c6de21c72dc6c1e1dae89a"  # was 61c4f18193adac7d146bc75c0f680430
$ws.Range("B553").Value = "58d85ba2051dd71507a5e4255d2e5b94"  # was 8317bc5e1079993b6d686cc7d773b4ef
$ws.Range("B572").Value = "f1eff8d1240251c266d684e4cbc1fca7"  # was 5ed55f8b2ae0bd9cea467720286f267b
$ws.Range("B584").Value = "a576e1b2662d1a21d6c1d37626fd4452"  # was e375d004872e7eac94fce210d9414135
$ws.Range("B666").Value = "6a504f8d367e29df8fe91b6e061f2350"  # was d0198b482e7ad0701fea272aba6657a8
$ws.Range("B729").Value = "27ed38bf1fbffac7273df8279ccba7ca"  # was b4db0bd5cfe9f51ea71702c7935a8b82
$ws.Range("B768").Value = "8a866f38cea4d509d812189b47eef642"  # was 856d009b685edcaa25e7aebd1e4cb92c
$ws.Range("B811").Value = "dbd952bba9bedbb15ced3d14a76bc9b0"  # was 5f1e48ea2ee37ac4a0cd6534daf28e1d
$ws.Range("B815").Value = "bd5b9380588c9dc7c9ba8123dc3cab76"  # was deeeabb02d47e448e34e5d3bbaeb8dad
$ws.Range("B816").Value = "1951623ae9020a139ec3467817acc2ab"  # was 831b12f239db1883cfb6a62cd480eabe
$ws.Range("B825").Value = "76fb08e3968f1341beee8c4d704ab1a6"  # was e0b748b7abab51601ff88878e1646e1d
$ws.Range("B827").Value = "fe391b223dd9b3e7fc6a5f6ebd9890a3"  # was e72e4ad52475855fd285dd2b5bbecbd4
$ws.Range("B874").Value = "d878f735a89572d2273c1e98708e28dd"  # was c9c849f03081bb7a17b5eba5feebb7ea
